$p = $ppt.ActivePresentation

function Fix-Title($slideIndex) {
    $s = $p.Slides.Item($slideIndex)
    $sh = $s.Shapes.Item(1)
    $tr = $sh.TextFrame.TextRange
    $lead = "Solve the following "

    # Drop the separately-styled "Solve the following " lead-in run (keeps the
    # remaining run's own rPr -- plain, no Roboto override) and rename the
    # slide title to the corrected text.
    $full = $tr.Text
    if ($full.IndexOf($lead) -eq 0) {
        $c1 = $tr.Characters(1, $lead.Length)
        $c1.Text = ""
    }
    $tr.Text = "Defining and Using Variables"
    Write-Host "Slide $slideIndex title -> $($tr.Text)"
}

# Slide 14 & 15: title "Solve the following algebra problem:" -> "Defining and Using Variables"
Fix-Title 14

# Slide 14: body paragraph - split "We need to be told what they mean (they need to be
# defined) first." into two runs: a red first clause, and a restyled/reworded second clause.
$s14 = $p.Slides.Item(14)
$body = $s14.Shapes.Item(5)
$tr = $body.TextFrame.TextRange

$oldSentence = "We need to be told what they mean (they need to be defined) first."
$full = $tr.Text
$idx = $full.IndexOf($oldSentence)
if ($idx -ge 0) {
    $startPos = $idx + 1

    $partA = "We need to be told what they mean"
    $cA = $tr.Characters($startPos, $partA.Length)
    $cA.Font.Color.RGB = 255

    $oldPartB = " (they need to be defined) first."
    $fullAfter = $tr.Text
    $idxB = $fullAfter.IndexOf($oldPartB)
    $cB = $tr.Characters($idxB + 1, $oldPartB.Length)
    $cB.Text = " (they need to be defined) BEFORE we can use them."
}
Write-Host "Slide 14 body -> $($tr.Text)"

Fix-Title 15
